# fixes in output dict creation, improved annotations, protocol changes
#
# The "strength" column was recomputed using an RMS-based measure instead
# of the previous raw measure, so the header label and the downstream
# numeric values (reactionTime, difference, strength columns) are updated
# to reflect the refreshed analysis output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label for column E (strength metric changed from raw to RMS)
$ws.Range("E1").Value = "strength (RMS)"

# Update recomputed data values (columns B-E) per refreshed analysis output
$ws.Range("B2").Value = 16.4
$ws.Range("D2").Value = 13.2
$ws.Range("E2").Value = 51.8

$ws.Range("B3").Value = 14
$ws.Range("D3").Value = 13.6
$ws.Range("E3").Value = 62.2

$ws.Range("E4").Value = 62

$ws.Range("B5").Value = 16.8
$ws.Range("D5").Value = 13.6
$ws.Range("E5").Value = 49.8

$ws.Range("B6").Value = 13.6
$ws.Range("D6").Value = 14
$ws.Range("E6").Value = 67.40000000000001

$ws.Range("B7").Value = 14
$ws.Range("D7").Value = 14
$ws.Range("E7").Value = 53.5

$ws.Range("B8").Value = 14.8
$ws.Range("D8").Value = 19.6
$ws.Range("E8").Value = 56.6

$ws.Range("B9").Value = 14
$ws.Range("C9").Value = 35.5
$ws.Range("D9").Value = 21.5
$ws.Range("E9").Value = 66

$ws.Range("B10").Value = 15.2
$ws.Range("D10").Value = 15.2
$ws.Range("E10").Value = 54.4

$ws.Range("B11").Value = 14
$ws.Range("D11").Value = 16
$ws.Range("E11").Value = 66.8

$ws.Range("B12").Value = 13.6
$ws.Range("D12").Value = 16.4
$ws.Range("E12").Value = 67.2

$ws.Range("B13").Value = 14.4
$ws.Range("D13").Value = 15.6
$ws.Range("E13").Value = 61.4

$ws.Range("B14").Value = 14
$ws.Range("D14").Value = 16
$ws.Range("E14").Value = 67.8

$ws.Range("B15").Value = 13
$ws.Range("D15").Value = 17
$ws.Range("E15").Value = 68.75

$ws.Range("E16").Value = 63

$ws.Range("B17").Value = 14.4
$ws.Range("D17").Value = 15.6
$ws.Range("E17").Value = 59.8

$ws.Range("B18").Value = 14.5
$ws.Range("D18").Value = 16
$ws.Range("E18").Value = 60.25

$ws.Range("B19").Value = 17.6
$ws.Range("D19").Value = 12.4
$ws.Range("E19").Value = 52.4
